$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 44, pushing the remaining log/notes rows down
# (same styling carries onto the new blank rows, matching the existing A-column
# date formatting already applied to that block).
$ws.Rows("44:46").Insert()

# Fill in the new time-log entry that now occupies row 43.
$ws.Range("A43").Value = 43539
$ws.Range("B43").Value = 0.5
$ws.Range("D43").Value = "Week 7/Project: worked on web.xml changes and package changes for authentication."

# Leave the selection where the author left it.
[void]$ws.Range("B44").Select()
